# Add two atributes 'Materia' - rename/recode 'especialidad' (D) column to 'seccion'
# mapping instruments -> sections (Cuerda / Viento / Percusion / Canto / Guitarra / Piano / Dirección)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header D1 from "especialidad" to "seccion"
$ws.Range("D1").Value = "seccion"

# 2. Recode each row's D value (instrument -> section)
$map = @{
    "Dirección" = "Dirección"
    "Guitarra"  = "Guitarra"
    "Piano"     = "Piano"
    "Cello"     = "Cuerda"
    "Violín"    = "Cuerda"
    "Clarinete" = "Viento"
    "Oboe"      = "Viento"
    "Canto"     = "Canto"
    "Viola"     = "Cuerda"
    "VIolín"    = "Cuerda"
}

for ($r = 2; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $current = $cell.Value2
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}

# Row 24 (profesor_id 1023, Denise Avila) was specifically corrected to "Percusion"
$ws.Cells.Item(24, 4).Value = "Percusion"

# 3. Move the selection to D25 (matches author's saved cursor position)
$ws.Range("D25").Select()
